# Weekly update: add a new "Caramelo" Piña price point (Femacal de La Calera)
# for 2023-12-05 (Excel serial 45265), inserted above the existing history.
# This pushes the rest of the historical rows (old 1150:1242) down by two
# rows (to 1152:1244) and expands the used range to A1:T1244.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current first data row of this block.
$ws.Rows("1150:1151").Insert()

# New row 1150: grade "Primera"
$ws.Range("A1150").Value = 3
$ws.Range("B1150").Value = "Femacal de La Calera"
$ws.Range("C1150").Value = "Coquimbo"
$ws.Range("D1150").Value = 45265
$ws.Range("E1150").Value = 5
$ws.Range("F1150").Value = "Fruta"
$ws.Range("G1150").Value = 100108
$ws.Range("H1150").Value = "Tropicales y subtropicales"
$ws.Range("I1150").Value = 100108005
$ws.Range("J1150").Value = "Piña"
$ws.Range("K1150").Value = "Caramelo"
$ws.Range("L1150").Value = "Primera"
$ws.Range("M1150").Value = 108
$ws.Range("N1150").Value = 23000
$ws.Range("O1150").Value = 23000
$ws.Range("P1150").Value = 23000
$ws.Range("Q1150").Value = "`$/caja 12 unidades"
$ws.Range("R1150").Value = "Ecuador"
$ws.Range("S1150").Value = 1917
$ws.Range("T1150").Value = 12

# New row 1151: grade "Segunda"
$ws.Range("A1151").Value = 3
$ws.Range("B1151").Value = "Femacal de La Calera"
$ws.Range("C1151").Value = "Coquimbo"
$ws.Range("D1151").Value = 45265
$ws.Range("E1151").Value = 5
$ws.Range("F1151").Value = "Fruta"
$ws.Range("G1151").Value = 100108
$ws.Range("H1151").Value = "Tropicales y subtropicales"
$ws.Range("I1151").Value = 100108005
$ws.Range("J1151").Value = "Piña"
$ws.Range("K1151").Value = "Caramelo"
$ws.Range("L1151").Value = "Segunda"
$ws.Range("M1151").Value = 108
$ws.Range("N1151").Value = 23000
$ws.Range("O1151").Value = 23000
$ws.Range("P1151").Value = 23000
$ws.Range("Q1151").Value = "`$/caja 14 unidades"
$ws.Range("R1151").Value = "Ecuador"
$ws.Range("S1151").Value = 1643
$ws.Range("T1151").Value = 14
